$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 2218.3333
$ws.Range("J112").Value = 2795.1428
$ws.Range("L112").Value = 8385.428400000001
$ws.Range("N112").Value = -10601.4284
$ws.Range("H113").Value = 1975.1177
$ws.Range("I113").Value = 1788.5
$ws.Range("J113").Value = 2000
$ws.Range("K113").Value = 1788.5
$ws.Range("L113").Value = 2000
$ws.Range("M113").Value = 1465.5
$ws.Range("N113").Value = -8508
$ws.Range("H116").Value = 1751001.8
$ws.Range("I116").Value = 3665022.2
$ws.Range("J116").Value = 3417.652
$ws.Range("K116").Value = 3665022.2
$ws.Range("L116").Value = 3417.652
$ws.Range("M116").Value = -3661580.2
$ws.Range("N116").Value = -10301.652
$ws.Range("H125").Value = 2016.6666
$ws.Range("I125").Value = 2007
$ws.Range("K125").Value = 18063
$ws.Range("M125").Value = -15603

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 19612882
$ws.Range("J32").Value = 17958.666
$ws.Range("L32").Value = 17958.666
$ws.Range("N32").Value = -18532.666
$ws.Range("H132").Value = 2376.2307
$ws.Range("J132").Value = 2549.5
$ws.Range("L132").Value = 7648.5
$ws.Range("N132").Value = -12708.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 212.08333
$ws.Range("I94").Value = 218.64285
$ws.Range("J94").Value = 202.9
$ws.Range("K94").Value = 218.64285
$ws.Range("L94").Value = 202.9
$ws.Range("M94").Value = 232.35715
$ws.Range("N94").Value = -1104.9
$ws.Range("H134").Value = 2670.8096
$ws.Range("I134").Value = 2468.4119
$ws.Range("K134").Value = 7405.2357
$ws.Range("M134").Value = -4870.2357

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2308.7344
$ws.Range("I31").Value = 1582.9608
$ws.Range("J31").Value = 5156
$ws.Range("K31").Value = 1582.9608
$ws.Range("L31").Value = 5156
$ws.Range("M31").Value = -1287.9608
$ws.Range("N31").Value = -5746
$ws.Range("H34").Value = 2308.7344
$ws.Range("I34").Value = 1582.9608
$ws.Range("J34").Value = 5156
$ws.Range("K34").Value = 1582.9608
$ws.Range("L34").Value = 5156
$ws.Range("M34").Value = -1380.9608
$ws.Range("N34").Value = -5560

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H97").Value = 935.0714
$ws.Range("I97").Value = 470.14285
$ws.Range("K97").Value = 1410.42855
$ws.Range("M97").Value = -914.4285500000001
$ws.Range("H109").Value = 3725.8333
$ws.Range("I109").Value = 1592
$ws.Range("J109").Value = 4501.773
$ws.Range("K109").Value = 4776
$ws.Range("L109").Value = 13505.319
$ws.Range("M109").Value = -3736
$ws.Range("N109").Value = -15585.319

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3204.9473
$ws.Range("I80").Value = 2542.3076
$ws.Range("J80").Value = 4640.6665
$ws.Range("K80").Value = 2542.3076
$ws.Range("L80").Value = 4640.6665
$ws.Range("M80").Value = -1544.3076
$ws.Range("N80").Value = -6636.6665
$ws.Range("H83").Value = 3204.9473
$ws.Range("I83").Value = 2542.3076
$ws.Range("J83").Value = 4640.6665
$ws.Range("K83").Value = 12711.538
$ws.Range("L83").Value = 23203.3325
$ws.Range("M83").Value = -7719.538
$ws.Range("N83").Value = -33187.3325
$ws.Range("H107").Value = 291.95
$ws.Range("I107").Value = 236
$ws.Range("K107").Value = 236
$ws.Range("M107").Value = 1684
$ws.Range("H126").Value = 2498.9333
$ws.Range("I126").Value = 1432.4
$ws.Range("K126").Value = 4297.200000000001
$ws.Range("M126").Value = -1827.200000000001
$ws.Range("H132").Value = 3504.6843
$ws.Range("I132").Value = 2908.5
$ws.Range("K132").Value = 8725.5
$ws.Range("M132").Value = -6195.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1898.6666
$ws.Range("I7").Value = 2042.0769
$ws.Range("J7").Value = 1665.625
$ws.Range("K7").Value = 2042.0769
$ws.Range("L7").Value = 1665.625
$ws.Range("M7").Value = -1930.0769
$ws.Range("N7").Value = -1889.625
$ws.Range("H40").Value = 4461
$ws.Range("I40").Value = 4556.933
$ws.Range("J40").Value = 4301.1113
$ws.Range("K40").Value = 4556.933
$ws.Range("L40").Value = 4301.1113
$ws.Range("M40").Value = -4420.933
$ws.Range("N40").Value = -4573.1113
$ws.Range("H61").Value = 4003
$ws.Range("I61").Value = 4937.4814
$ws.Range("K61").Value = 4937.4814
$ws.Range("M61").Value = -4735.4814
$ws.Range("H82").Value = 1635
$ws.Range("I82").Value = 1473.125
$ws.Range("J82").Value = 2066.6667
$ws.Range("K82").Value = 1473.125
$ws.Range("L82").Value = 2066.6667
$ws.Range("M82").Value = -1112.125
$ws.Range("N82").Value = -2788.6667
$ws.Range("H85").Value = 1635
$ws.Range("I85").Value = 1473.125
$ws.Range("J85").Value = 2066.6667
$ws.Range("K85").Value = 1473.125
$ws.Range("L85").Value = 2066.6667
$ws.Range("M85").Value = -225.125
$ws.Range("N85").Value = -4562.6667
$ws.Range("H113").Value = 4003
$ws.Range("I113").Value = 4937.4814
$ws.Range("K113").Value = 4937.4814
$ws.Range("M113").Value = -2767.4814
$ws.Range("H126").Value = 1898.6666
$ws.Range("I126").Value = 2042.0769
$ws.Range("J126").Value = 1665.625
$ws.Range("K126").Value = 6126.2307
$ws.Range("L126").Value = 4996.875
$ws.Range("M126").Value = -3656.2307
$ws.Range("N126").Value = -9936.875
$ws.Range("H132").Value = 4015.4614
$ws.Range("I132").Value = 3522.5557
$ws.Range("K132").Value = 10567.6671
$ws.Range("M132").Value = -8037.667099999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 436.4375
$ws.Range("I113").Value = 306.3846
$ws.Range("K113").Value = 919.1537999999999
$ws.Range("M113").Value = 1250.8462
$ws.Range("H132").Value = 1995.1562
$ws.Range("I132").Value = 1260.1052
$ws.Range("J132").Value = 3069.4614
$ws.Range("K132").Value = 3780.3156
$ws.Range("L132").Value = 9208.3842
$ws.Range("M132").Value = -1250.3156
$ws.Range("N132").Value = -14268.3842
$ws.Range("H135").Value = 32000
$ws.Range("J135").Value = 32000
$ws.Range("L135").Value = 32000
$ws.Range("N135").Value = -42140
$ws.Range("H136").Value = 1510.3864
$ws.Range("I136").Value = 1269.1562
$ws.Range("J136").Value = 2153.6667
$ws.Range("K136").Value = 3807.4686
$ws.Range("L136").Value = 6461.000100000001
$ws.Range("M136").Value = -1257.4686
$ws.Range("N136").Value = -11561.0001
